# Apply the commit's changes to the "2€" sheet of the Latvia commemorative
# UNC varieties workbook:
#   - Mark the H3 variety count as 1 (was 0)
#   - Mark the I16, I17, I18 variety counts as 1 (was 0)
#   - Move/leave the active selection on the frozen bottom-right pane at H28

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")
$ws.Activate()

$ws.Range("H3").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("I18").Value = 1

# Reproduce the new selected cell recorded in the saved view (bottom-right
# frozen pane, since the sheet is split/frozen at column L / row 3).
$ws.Range("H28").Select()

# Best-effort: keep window geometry in sync with the commit (not all hosts
# persist these window-chrome attributes, but set them for parity).
$excel.ActiveWindow.Left = -110
$excel.ActiveWindow.Top = -110
$excel.ActiveWindow.Width = 38620
$excel.ActiveWindow.Height = 21220
